$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Relatório Compra 16-03-2025"

$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 44
$ws.Range("F2").Value = 88
$ws.Range("G2").Value = 1.02
$ws.Range("H2").Value = 1.62
$ws.Range("J2").Value = 90
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "BANANA NANICA"
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 40
$ws.Range("G3").Value = 2.5
$ws.Range("H3").Value = 3.95
$ws.Range("J3").Value = 100
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "BANANA PRATA"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 45
$ws.Range("F4").Value = 20
$ws.Range("G4").Value = 2.25
$ws.Range("H4").Value = 3.55
$ws.Range("I4").Value = 0.5788
$ws.Range("J4").Value = 45
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "MAMAO PAPAIA"
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 30
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 3.16
$ws.Range("I5").Value = 0.58
$ws.Range("J5").Value = 60
$ws.Range("A6").Value = 53
$ws.Range("B6").Value = "BATATA LAVADA"
$ws.Range("D6").Value = 90
$ws.Range("E6").Value = 25
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = 3.6
$ws.Range("H6").Value = 5.87
$ws.Range("I6").Value = 0.63
$ws.Range("J6").Value = 180
$ws.Range("A7").Value = 13
$ws.Range("B7").Value = "CAQUI CHOCOLATE"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 80
$ws.Range("F7").Value = 20
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = 13.04
$ws.Range("I7").Value = 0.63
$ws.Range("J7").Value = 160
$ws.Range("A8").Value = 28
$ws.Range("B8").Value = "MELANCIA"
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 10
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = 48.9
$ws.Range("I8").Value = 0.63
$ws.Range("J8").Value = 300
$ws.Range("A9").Value = 62
$ws.Range("B9").Value = "CHUCHU"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 45
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 20
$ws.Range("G9").Value = 4.5
$ws.Range("H9").Value = 7.33
$ws.Range("I9").Value = 0.63
$ws.Range("J9").Value = 90
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "LARANJA PERA"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 50
$ws.Range("E10").Value = 22
$ws.Range("F10").Value = 44
$ws.Range("G10").Value = 2.27
$ws.Range("H10").Value = 3.7
$ws.Range("I10").Value = 0.63
$ws.Range("J10").Value = 100
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "MACA FUJI"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 90
$ws.Range("E11").Value = 20
$ws.Range("F11").Value = 40
$ws.Range("G11").Value = 4.5
$ws.Range("H11").Value = 8.1
$ws.Range("I11").Value = 0.8
$ws.Range("J11").Value = 180
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "ALFACE CRESPA"
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = 18
$ws.Range("F12").Value = 54
$ws.Range("G12").Value = 1.67
$ws.Range("H12").Value = 2.63
$ws.Range("I12").Value = 0.58
$ws.Range("J12").Value = 90
